$d = $word.ActiveDocument

# paragraph 12
$p12 = $d.Paragraphs.Item(12)
$r12 = $p12.Range
$r12.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 16 commits</w:t></w:r><w:r><w:t xml:space="preserve"> as of Apr 15, 2018.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# paragraph 17
$p17 = $d.Paragraphs.Item(17)
$r17 = $p17.Range
$r17.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D97119" w:rsidRDefault="00D97119"><w:r><w:t xml:space="preserve">Created </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Youtube</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> account.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# paragraph 19
$p19 = $d.Paragraphs.Item(19)
$r19 = $p19.Range
$r19.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 14</w:t></w:r><w:r><w:t xml:space="preserve"> commits</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>as of Apr 15, 2018.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# paragraph 26
$p26 = $d.Paragraphs.Item(26)
$r26 = $p26.Range
$r26.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 9</w:t></w:r><w:r><w:t xml:space="preserve"> commits</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>as of Apr 15, 2018.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# paragraph 32
$p32 = $d.Paragraphs.Item(32)
$r32 = $p32.Range
$r32.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 3</w:t></w:r><w:r><w:t xml:space="preserve"> commits</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>as of Apr 15, 2018.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# paragraph 37
$p37 = $d.Paragraphs.Item(37)
$r37 = $p37.Range
$r37.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00881DC4" w:rsidRDefault="00881DC4"><w:r><w:t xml:space="preserve">Created the group </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Trello</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> account. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# paragraph 39
$p39 = $d.Paragraphs.Item(39)
$r39 = $p39.Range
$r39.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00BE394D" w:rsidRDefault="00BE394D"><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 3</w:t></w:r><w:r><w:t xml:space="preserve"> commits</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>as of Apr 15, 2018.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

